$d = $word.ActiveDocument

$replacements = @(
    @("2024-05-18 Saturday", "2024-05-19 Sunday"),
    @("825÷4=206, 1", "894÷6=149, 0"),
    @("469÷7=67, 0", "390÷4=97, 2"),
    @("291÷5=58, 1", "112÷2=56, 0"),
    @("656÷8=82, 0", "738÷2=369, 0"),
    @("776÷7=110, 6", "438÷4=109, 2"),
    @("252÷9=28, 0", "739÷5=147, 4"),
    @("730÷8=91, 2", "376÷6=62, 4"),
    @("824÷4=206, 0", "720÷7=102, 6"),
    @("185÷3=61, 2", "313÷2=156, 1"),
    @("630÷7=90, 0", "974÷4=243, 2"),
    @("960÷4=240, 0", "338÷4=84, 2"),
    @("696÷6=116, 0", "637÷3=212, 1"),
    @("374÷7=53, 3", "913÷6=152, 1"),
    @("357÷6=59, 3", "315÷7=45, 0"),
    @("878÷6=146, 2", "770÷9=85, 5"),
    @("551÷6=91, 5", "683÷3=227, 2"),
    @("883÷7=126, 1", "879÷5=175, 4"),
    @("573÷2=286, 1", "526÷9=58, 4"),
    @("412÷3=137, 1", "337÷6=56, 1"),
    @("819÷3=273, 0", "400÷9=44, 4"),
    @("354÷2=177, 0", "848÷7=121, 1"),
    @("442÷2=221, 0", "801÷4=200, 1"),
    @("335÷8=41, 7", "222÷4=55, 2"),
    @("524÷3=174, 2", "141÷5=28, 1"),
    @("164÷6=27, 2", "251÷2=125, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done applying $($replacements.Count) replacements"
